$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39
$ws.Range("A39").Value = 7
$ws.Range("C39").Value = "extr_low_pass_filter_cond1"
$ws.Range("D39").Value = 35

# Row 40
$ws.Range("A40").Value = 7
$ws.Range("C40").Value = "extr_high_pass_filter_cond1"
$ws.Range("D40").Value = 45

# Row 41
$ws.Range("A41").Value = 7
$ws.Range("C41").Value = "extr_low_pass_filter_cond2"
$ws.Range("D41").Value = 22

# Row 42
$ws.Range("A42").Value = 7
$ws.Range("C42").Value = "extr_high_pass_filter_cond2"
$ws.Range("D42").Value = 32

# Row 43
$ws.Range("A43").Value = 3
$ws.Range("C43").Value = "ICA_eye"
$ws.Range("D43").Value = 0.8

# Row 44
$ws.Range("A44").Value = 3
$ws.Range("C44").Value = "ICA_brain"
$ws.Range("D44").Value = 0.1

# Apply style to C39:C47 like C38 (style index 1)
$ws.Range("C39:C47").Style = $ws.Range("C38").Style

# Selection and view
$ws.Range("D48").Select
